$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 90, shifting existing rows 90-137 down to 91-138.
$ws.Rows.Item(90).Insert()

# Populate the newly inserted row 90 with the new record's data.
$ws.Cells.Item(90, 1).Value = 11
$ws.Cells.Item(90, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(90, 3).Value = "Bíobío"
$ws.Cells.Item(90, 4).Value = 44846
$ws.Cells.Item(90, 5).Value = 8
$ws.Cells.Item(90, 6).Value = 100112021
$ws.Cells.Item(90, 7).Value = "Ají"
$ws.Cells.Item(90, 8).Value = "Inferno"
$ws.Cells.Item(90, 9).Value = "Primera"
$ws.Cells.Item(90, 10).Value = 35
$ws.Cells.Item(90, 11).Value = 18000
$ws.Cells.Item(90, 12).Value = 19000
$ws.Cells.Item(90, 13).Value = 18429
$ws.Cells.Item(90, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(90, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(90, 16).Value = 1843
$ws.Cells.Item(90, 17).Value = 10
$ws.Cells.Item(90, 18).Value = "Hortaliza"
